$wb = $excel.ActiveWorkbook

# --- Sheet "FAST v8.09.00" (first sheet) ---
$ws1 = $wb.Worksheets.Item("FAST v8.09.00")

# Updated "as-of" date for the registry (B2)
$ws1.Range("B2").Value = 41912

# New dependency source-file entries added in this order so that the
# shared-string table gets populated in the same sequence as the
# original authoring (row 25, then 15, then 16, then 14).
$ws1.Range("C25").Value = "v2.04.00a-bjj"
$ws1.Range("D25").Value = 41912
$ws1.Range("G25").Value = 261
$ws1.Range("H25").Value = 41908

$ws1.Range("C15").Value = "v1.01.07a-bjj"
$ws1.Range("D15").Value = 41912

$ws1.Range("C16").Value = "v14.03.00a-bjj"
$ws1.Range("D16").Value = 41912

$ws1.Range("C14").Value = "v8.09.00a-bjj"
$ws1.Range("D14").Value = 41912

# I17 becomes a "not applicable" hatched cell (keeps the date number
# format, just no value, with the gray125 pattern fill).
$ws1.Range("I17").Interior.Pattern = 17

# --- Sheet "FAST v8.08.00" (second sheet) ---
$ws2 = $wb.Worksheets.Item("FAST v8.08.00")
$ws2.Range("C25").Select() | Out-Null

# Update the active selection left by the last editor; re-select on
# sheet 1 last so it remains the active/tabSelected sheet.
$ws1.Range("I17").Select() | Out-Null
